# Update the table style ("Table Design" gallery selection) on the table
# that currently uses the old custom table style GUID, switching it to the
# new built-in PowerPoint table style GUID.
#
# {EACEACA1-F5E5-4F06-AAA9-B033EB3F6054}  ->  {A7E4DC3A-D4D8-4C63-BD50-47A2D0A89B39}

$OldStyleId = "{EACEACA1-F5E5-4F06-AAA9-B033EB3F6054}"
$NewStyleId = "{A7E4DC3A-D4D8-4C63-BD50-47A2D0A89B39}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $OldStyleId) {
                $table.ApplyStyle($NewStyleId)
            }
        }
    }
}
